# first sweep cleaning data columns to conform to specs--done by chase
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the space-separated values in columns D (experimentDesign) and
# F (strain) with underscore-separated equivalents, across every data row
# in the sheet's used range. Matching on the original text means only the
# cells that actually still hold the old values get touched.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    if ($dCell.Value2 -eq "Environmental Perturbation") {
        $dCell.Value2 = "Environmental_Perturbation"
    }

    $fCell = $ws.Cells.Item($r, 6)
    if ($fCell.Value2 -eq "KN99 alpha") {
        $fCell.Value2 = "KN99_alpha"
    }
}

# Update the sheet's view/selection state: scroll so column F is the
# leftmost visible column, and move the active selection to O15.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 6
$ws.Range("O15").Select()
